$d = $word.ActiveDocument

$range = $d.Content
$found = $range.Find.Execute("{ {0} }", $true, $false, $false, $false, $false, $true, 1, $false, "{ {0} };", 2)
